# issue #5: property land done
#
# Normalizes free-text date / address / land-lot strings (drops stray
# internal spaces and, for land-lot numbers, the "-" separator), retitles
# the 土地 (land) sheet's header row from human-readable Chinese labels to
# the canonical schema field names, and appends the standard metadata
# columns (property_category/category/date/legislator_name/legislator_id/
# source_file/index) to the 土地 sheet, mirroring the 股票 sheet layout.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 土地 (land) sheet
# ---------------------------------------------------------------------
$land = $wb.Worksheets.Item("土地")

# Re-header row 1 with schema field names instead of Chinese captions.
$land.Range("B1").Value = "name"
$land.Range("C1").Value = "area"
$land.Range("D1").Value = "share_portion"
$land.Range("E1").Value = "owner"
$land.Range("F1").Value = "register_date"
$land.Range("G1").Value = "register_reason"
$land.Range("H1").Value = "acquire_value"

# Strip the "-" out of the land-lot numbers (and the trailing space before
# 地號).
$land.Range("B2").Value = "臺東縣東河鄉都蘭段10420000地號"
$land.Range("B3").Value = "臺東縣成功鎮忠智段11700001地號"
$land.Range("B4").Value = "臺東縣台東市豐康段07330119地號"
$land.Range("B5").Value = "臺東縣台東市新生段00020006地號"

# Collapse the stray space in the register-date strings.
$land.Range("F2").Value = "89年09月27日"
$land.Range("F3").Value = "86年11月03日"
$land.Range("F4").Value = "90年01月12日"
$land.Range("F5").Value = "94年04月15日"

# Append the standard metadata columns (same shape as the 股票 sheet).
$landRows = @(13, 14, 15, 16)
for ($i = 0; $i -lt $landRows.Length; $i++) {
    $r = $i + 2
    $land.Cells.Item($r, 9).Value = "land"
    $land.Cells.Item($r, 10).Value = "normal"
    $land.Cells.Item($r, 11).Value = "2012-03-06"
    $land.Cells.Item($r, 12).Value = "廖國棟"
    $land.Cells.Item($r, 13).Value = 962
    $land.Cells.Item($r, 14).Value = "tmpec731"
    $land.Cells.Item($r, 15).Value = $landRows[$i]
}

$land.Range("I1").Value = "property_category"
$land.Range("J1").Value = "category"
$land.Range("K1").Value = "date"
$land.Range("L1").Value = "legislator_name"
$land.Range("M1").Value = "legislator_id"
$land.Range("N1").Value = "source_file"
$land.Range("O1").Value = "index"

# ---------------------------------------------------------------------
# 建物 (building) sheet
# ---------------------------------------------------------------------
$building = $wb.Worksheets.Item("建物")
$building.Range("F2").Value = "86年11月03日"
$building.Range("F3").Value = "94年09月15日"
$building.Range("F4").Value = "90年01月12F1"

# ---------------------------------------------------------------------
# 汽車 (vehicle) sheet
# ---------------------------------------------------------------------
$vehicle = $wb.Worksheets.Item("汽車")
$vehicle.Range("E2").Value = "99年09月01日"
$vehicle.Range("E3").Value = "93年01月06H"

# ---------------------------------------------------------------------
# 債務 (debt) sheet
# ---------------------------------------------------------------------
$debt = $wb.Worksheets.Item("債務")
$debt.Range("D2").Value = "土地銀行台東分行臺東縣台東市中華路"
$debt.Range("D3").Value = "華南銀行台東分行臺東縣台東市中華路"
$debt.Range("D4").Value = "台灣企銀台東分行臺東縣台東市中華路"
$debt.Range("D5").Value = "台新銀行三重分行新北市三重區正義北路"
$debt.Range("D6").Value = "台東縣都蘭農會臺東縣東河鄉都蘭村都蘭"
$debt.Range("D7").Value = "土地銀行台東分行臺東縣台東市中華路"
$debt.Range("D8").Value = "台東縣都蘭農會臺東縣東河鄉都蘭村都蘭"

$debt.Range("F2").Value = "98年05月06日"
$debt.Range("F3").Value = "90年01月12H"
$debt.Range("F4").Value = "94年10月11曰"
$debt.Range("F5").Value = "99年08月27日"
$debt.Range("F6").Value = "93年08月31日"
$debt.Range("F7").Value = "100年07月01日"
$debt.Range("F8").Value = "96年03月13曰"
